$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 2 new columns before column D, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# --- Copy formatting for the new D:E columns in rows 7:35 ---
$ws.Range("F7:G7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("F8:G35").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- Copy formatting for the new D:E columns in rows 38:77 ---
$ws.Range("F38:G38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("F39:G77").Copy()
$ws.Range("D39").PasteSpecial(-4122)

# --- Copy formatting for the new D:E columns in rows 80:102 ---
$ws.Range("F80:G80").Copy()
$ws.Range("D80").PasteSpecial(-4122)
$ws.Range("F81:G102").Copy()
$ws.Range("D81").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Set new quarter data values in columns D and E ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 542800
$ws.Range("E8").Value = 536900
$ws.Range("D9").Value = 88500
$ws.Range("E9").Value = 83300
$ws.Range("D10").Value = 454300
$ws.Range("E10").Value = 453600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 17800
$ws.Range("E14").Value = 204100
$ws.Range("D15").Value = 174100
$ws.Range("E15").Value = 173400
$ws.Range("D17").Value = 328100
$ws.Range("E17").Value = 503700
$ws.Range("D18").Value = 214700
$ws.Range("E18").Value = 33200
$ws.Range("D20").Value = -324900
$ws.Range("E20").Value = -299000
$ws.Range("D21").Value = 63900
$ws.Range("E21").Value = -92300
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -110200
$ws.Range("E23").Value = -265800
$ws.Range("D24").Value = 200
$ws.Range("E24").Value = 107900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -110400
$ws.Range("E26").Value = -373600
$ws.Range("D27").Value = -111300
$ws.Range("E27").Value = -374600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 324900
$ws.Range("E32").Value = 299000
$ws.Range("D33").Value = -111300
$ws.Range("E33").Value = -374600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -111300
$ws.Range("E35").Value = -374600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 485100
$ws.Range("E41").Value = 660800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 316400
$ws.Range("E43").Value = 291100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 46100
$ws.Range("E45").Value = 52300
$ws.Range("D46").Value = 847700
$ws.Range("E46").Value = 1004200
$ws.Range("D47").Value = 96100
$ws.Range("E47").Value = 86000
$ws.Range("D48").Value = 5511700
$ws.Range("E48").Value = 5608800
$ws.Range("D49").Value = 5384600
$ws.Range("E49").Value = 5394300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 401400
$ws.Range("E52").Value = 421100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 12241500
$ws.Range("E54").Value = 12514300
$ws.Range("D57").Value = 108100
$ws.Range("E57").Value = 96300
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 552100
$ws.Range("E59").Value = 468800
$ws.Range("D60").Value = 660200
$ws.Range("E60").Value = 565200
$ws.Range("D61").Value = 14028400
$ws.Range("E61").Value = 14270500
$ws.Range("D62").Value = 1635600
$ws.Range("E62").Value = 1654300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 16338500
$ws.Range("E66").Value = 16505600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -6606400
$ws.Range("E72").Value = -6495100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -4097000
$ws.Range("E76").Value = -3991200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -111300
$ws.Range("E81").Value = -374600
$ws.Range("D83").Value = 174100
$ws.Range("E83").Value = 173400
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 170000
$ws.Range("E89").Value = 36100
$ws.Range("D91").Value = -79700
$ws.Range("E91").Value = -55600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -66200
$ws.Range("E94").Value = -74800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -278400
$ws.Range("E100").Value = 258700
$ws.Range("D101").Value = 500
$ws.Range("E101").Value = -1700
$ws.Range("D102").Value = -174100
$ws.Range("E102").Value = 218300
